$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 124,3
$arr[0,0] = 1
$arr[0,1] = 1
$arr[0,2] = 0
$arr[1,0] = 1
$arr[1,1] = 2
$arr[1,2] = 83
$arr[2,0] = 1
$arr[2,1] = 3
$arr[2,2] = 84
$arr[3,0] = 1
$arr[3,1] = 4
$arr[3,2] = 89
$arr[4,0] = 1
$arr[4,1] = 5
$arr[4,2] = 86
$arr[5,0] = 1
$arr[5,1] = 6
$arr[5,2] = 96
$arr[6,0] = 1
$arr[6,1] = 7
$arr[6,2] = 104
$arr[7,0] = 1
$arr[7,1] = 8
$arr[7,2] = 109
$arr[8,0] = 1
$arr[8,1] = 9
$arr[8,2] = 113
$arr[9,0] = 1
$arr[9,1] = 10
$arr[9,2] = 114
$arr[10,0] = 1
$arr[10,1] = 11
$arr[10,2] = 115
$arr[11,0] = 1
$arr[11,1] = 12
$arr[11,2] = 118
$arr[12,0] = 1
$arr[12,1] = 13
$arr[12,2] = 117
$arr[13,0] = 1
$arr[13,1] = 14
$arr[13,2] = 120
$arr[14,0] = 1
$arr[14,1] = 15
$arr[14,2] = 119
$arr[15,0] = 1
$arr[15,1] = 16
$arr[15,2] = 116
$arr[16,0] = 1
$arr[16,1] = 17
$arr[16,2] = 108
$arr[17,0] = 1
$arr[17,1] = 18
$arr[17,2] = 110
$arr[18,0] = 1
$arr[18,1] = 19
$arr[18,2] = 107
$arr[19,0] = 1
$arr[19,1] = 20
$arr[19,2] = 100
$arr[20,0] = 1
$arr[20,1] = 21
$arr[20,2] = 99
$arr[21,0] = 1
$arr[21,1] = 22
$arr[21,2] = 103
$arr[22,0] = 1
$arr[22,1] = 23
$arr[22,2] = 77
$arr[23,0] = 1
$arr[23,1] = 24
$arr[23,2] = 74
$arr[24,0] = 1
$arr[24,1] = 25
$arr[24,2] = 82
$arr[25,0] = 1
$arr[25,1] = 26
$arr[25,2] = 71
$arr[26,0] = 1
$arr[26,1] = 27
$arr[26,2] = 111
$arr[27,0] = 1
$arr[27,1] = 28
$arr[27,2] = 30
$arr[28,0] = 1
$arr[28,1] = 29
$arr[28,2] = 112
$arr[29,0] = 2
$arr[29,1] = 1
$arr[29,2] = 0
$arr[30,0] = 2
$arr[30,1] = 2
$arr[30,2] = 85
$arr[31,0] = 2
$arr[31,1] = 3
$arr[31,2] = 91
$arr[32,0] = 2
$arr[32,1] = 4
$arr[32,2] = 105
$arr[33,0] = 2
$arr[33,1] = 5
$arr[33,2] = 106
$arr[34,0] = 2
$arr[34,1] = 6
$arr[34,2] = 98
$arr[35,0] = 2
$arr[35,1] = 7
$arr[35,2] = 102
$arr[36,0] = 2
$arr[36,1] = 8
$arr[36,2] = 101
$arr[37,0] = 2
$arr[37,1] = 9
$arr[37,2] = 95
$arr[38,0] = 2
$arr[38,1] = 10
$arr[38,2] = 88
$arr[39,0] = 2
$arr[39,1] = 11
$arr[39,2] = 87
$arr[40,0] = 2
$arr[40,1] = 12
$arr[40,2] = 93
$arr[41,0] = 2
$arr[41,1] = 13
$arr[41,2] = 94
$arr[42,0] = 2
$arr[42,1] = 14
$arr[42,2] = 97
$arr[43,0] = 2
$arr[43,1] = 15
$arr[43,2] = 92
$arr[44,0] = 2
$arr[44,1] = 16
$arr[44,2] = 90
$arr[45,0] = 2
$arr[45,1] = 17
$arr[45,2] = 76
$arr[46,0] = 2
$arr[46,1] = 18
$arr[46,2] = 75
$arr[47,0] = 2
$arr[47,1] = 19
$arr[47,2] = 78
$arr[48,0] = 2
$arr[48,1] = 20
$arr[48,2] = 80
$arr[49,0] = 2
$arr[49,1] = 21
$arr[49,2] = 70
$arr[50,0] = 2
$arr[50,1] = 22
$arr[50,2] = 62
$arr[51,0] = 2
$arr[51,1] = 23
$arr[51,2] = 66
$arr[52,0] = 2
$arr[52,1] = 24
$arr[52,2] = 63
$arr[53,0] = 2
$arr[53,1] = 25
$arr[53,2] = 64
$arr[54,0] = 2
$arr[54,1] = 26
$arr[54,2] = 79
$arr[55,0] = 2
$arr[55,1] = 27
$arr[55,2] = 56
$arr[56,0] = 2
$arr[56,1] = 28
$arr[56,2] = 51
$arr[57,0] = 2
$arr[57,1] = 29
$arr[57,2] = 49
$arr[58,0] = 2
$arr[58,1] = 30
$arr[58,2] = 42
$arr[59,0] = 2
$arr[59,1] = 31
$arr[59,2] = 39
$arr[60,0] = 2
$arr[60,1] = 32
$arr[60,2] = 47
$arr[61,0] = 2
$arr[61,1] = 33
$arr[61,2] = 48
$arr[62,0] = 2
$arr[62,1] = 34
$arr[62,2] = 40
$arr[63,0] = 2
$arr[63,1] = 35
$arr[63,2] = 38
$arr[64,0] = 2
$arr[64,1] = 36
$arr[64,2] = 33
$arr[65,0] = 2
$arr[65,1] = 37
$arr[65,2] = 14
$arr[66,0] = 2
$arr[66,1] = 38
$arr[66,2] = 6
$arr[67,0] = 2
$arr[67,1] = 39
$arr[67,2] = 5
$arr[68,0] = 2
$arr[68,1] = 40
$arr[68,2] = 2
$arr[69,0] = 2
$arr[69,1] = 41
$arr[69,2] = 21
$arr[70,0] = 2
$arr[70,1] = 42
$arr[70,2] = 18
$arr[71,0] = 3
$arr[71,1] = 1
$arr[71,2] = 0
$arr[72,0] = 3
$arr[72,1] = 2
$arr[72,2] = 73
$arr[73,0] = 3
$arr[73,1] = 3
$arr[73,2] = 69
$arr[74,0] = 3
$arr[74,1] = 4
$arr[74,2] = 61
$arr[75,0] = 3
$arr[75,1] = 5
$arr[75,2] = 55
$arr[76,0] = 3
$arr[76,1] = 6
$arr[76,2] = 52
$arr[77,0] = 3
$arr[77,1] = 7
$arr[77,2] = 32
$arr[78,0] = 3
$arr[78,1] = 8
$arr[78,2] = 35
$arr[79,0] = 3
$arr[79,1] = 9
$arr[79,2] = 41
$arr[80,0] = 3
$arr[80,1] = 10
$arr[80,2] = 43
$arr[81,0] = 3
$arr[81,1] = 11
$arr[81,2] = 28
$arr[82,0] = 3
$arr[82,1] = 12
$arr[82,2] = 24
$arr[83,0] = 3
$arr[83,1] = 13
$arr[83,2] = 22
$arr[84,0] = 3
$arr[84,1] = 14
$arr[84,2] = 17
$arr[85,0] = 3
$arr[85,1] = 15
$arr[85,2] = 15
$arr[86,0] = 3
$arr[86,1] = 16
$arr[86,2] = 12
$arr[87,0] = 3
$arr[87,1] = 17
$arr[87,2] = 11
$arr[88,0] = 3
$arr[88,1] = 18
$arr[88,2] = 9
$arr[89,0] = 3
$arr[89,1] = 19
$arr[89,2] = 4
$arr[90,0] = 3
$arr[90,1] = 20
$arr[90,2] = 3
$arr[91,0] = 3
$arr[91,1] = 21
$arr[91,2] = 8
$arr[92,0] = 3
$arr[92,1] = 22
$arr[92,2] = 10
$arr[93,0] = 3
$arr[93,1] = 23
$arr[93,2] = 13
$arr[94,0] = 3
$arr[94,1] = 24
$arr[94,2] = 25
$arr[95,0] = 3
$arr[95,1] = 25
$arr[95,2] = 34
$arr[96,0] = 3
$arr[96,1] = 26
$arr[96,2] = 36
$arr[97,0] = 3
$arr[97,1] = 27
$arr[97,2] = 27
$arr[98,0] = 3
$arr[98,1] = 28
$arr[98,2] = 44
$arr[99,0] = 3
$arr[99,1] = 29
$arr[99,2] = 45
$arr[100,0] = 3
$arr[100,1] = 30
$arr[100,2] = 31
$arr[101,0] = 3
$arr[101,1] = 31
$arr[101,2] = 29
$arr[102,0] = 3
$arr[102,1] = 32
$arr[102,2] = 23
$arr[103,0] = 3
$arr[103,1] = 33
$arr[103,2] = 19
$arr[104,0] = 3
$arr[104,1] = 34
$arr[104,2] = 20
$arr[105,0] = 3
$arr[105,1] = 35
$arr[105,2] = 16
$arr[106,0] = 3
$arr[106,1] = 36
$arr[106,2] = 7
$arr[107,0] = 3
$arr[107,1] = 37
$arr[107,2] = 1
$arr[108,0] = 4
$arr[108,1] = 1
$arr[108,2] = 0
$arr[109,0] = 4
$arr[109,1] = 2
$arr[109,2] = 81
$arr[110,0] = 4
$arr[110,1] = 3
$arr[110,2] = 72
$arr[111,0] = 4
$arr[111,1] = 4
$arr[111,2] = 60
$arr[112,0] = 4
$arr[112,1] = 5
$arr[112,2] = 65
$arr[113,0] = 4
$arr[113,1] = 6
$arr[113,2] = 68
$arr[114,0] = 4
$arr[114,1] = 7
$arr[114,2] = 67
$arr[115,0] = 4
$arr[115,1] = 8
$arr[115,2] = 59
$arr[116,0] = 4
$arr[116,1] = 9
$arr[116,2] = 58
$arr[117,0] = 4
$arr[117,1] = 10
$arr[117,2] = 57
$arr[118,0] = 4
$arr[118,1] = 11
$arr[118,2] = 54
$arr[119,0] = 4
$arr[119,1] = 12
$arr[119,2] = 26
$arr[120,0] = 4
$arr[120,1] = 13
$arr[120,2] = 53
$arr[121,0] = 4
$arr[121,1] = 14
$arr[121,2] = 50
$arr[122,0] = 4
$arr[122,1] = 15
$arr[122,2] = 37
$arr[123,0] = 4
$arr[123,1] = 16
$arr[123,2] = 46
$ws.Range("A2:C125").Value = $arr
Write-Host "Updated A2:C125 with 124 rows"
